$d = $word.ActiveDocument

# --- Paragraph 1: split "{{#each chapters}}{{title}}" into proofed runs ---
$p1 = $d.Paragraphs(1).Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rStyle w:val="Heading1Char"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{{#each chapters</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}}</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Heading1Char"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rStyle w:val="Heading1Char"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rStyle w:val="Heading1Char"/></w:rPr><w:t>title</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rStyle w:val="Heading1Char"/></w:rPr><w:t>}}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p1.InsertXML($xml)

# Re-apply the "Heading 1 Char" run style lost on direct run-level rStyle during InsertXML
$rg = $d.Range(18, 27)
$rg.Style = "Heading 1 Char"

# --- Paragraph 3: drop the leftover _GoBack bookmark ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

Write-Output "ok"
